$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Analysis method" header to "analysis_method" for uniformity
$ws.Range("A1").Value = "analysis_method"

# New column M: option to export the complete dep object output
$ws.Range("M1").Value = "complete_output"

# Switch the example analysis method value from MaxQuant to DiaNN
$ws.Range("A2").Value = "DiaNN"

$ws.Range("M2").Value = $true

# Add a TRUE/FALSE list data validation on the new cell, matching the
# style of the workbook's other list validations
[void]$ws.Range("M2").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Update the active selection to the newly added cell
[void]$ws.Range("M2").Select()
